$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "parameter values" (sheet1)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("parameter values")

# Row 43 ("N2O reduced emissions in land managed with biochar") values change
$ws1.Range("B43").Value = 0.39
$ws1.Range("C43").Value = -0.02
$ws1.Range("D43").Value = -0.46

# Row 43 loses its special yellow-highlight/percent style and now matches the
# plain percent style used by the rest of the table (same as row 44 B:D).
$ws1.Range("B44:D44").Copy()
$ws1.Range("B43:D43").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 43 is shorter now (one citation instead of a long multi-part one) -
# update its row height to match.
$ws1.Rows.Item(43).RowHeight = 63.75

# The source for the N2O row now cites Cayuela et al. 2014 instead of the
# Borchard et al. 2019 meta-analysis.
$ws1.Range("F43").Value = "Cayuela, M.L., Van Zwieten, L., Singh, B.P., Jeffery, S., Roig, A. and Sánchez-Monedero, M.A., 2014. Biochar's role in mitigating soil nitrous oxide emissions: A review and meta-analysis. Agriculture, Ecosystems & Environment, 191, pp.5-16."

# Selection on this sheet moves from G43 to D43.
$ws1.Activate()
$ws1.Range("D43").Select()

# ---------------------------------------------------------------------------
# Sheet "scenario tracker" (sheet2)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("scenario tracker")

# Scenarios that were "waiting on soil N2O flux coefficients" (with or
# without "need to make R file") now have their R files created.
$ws2.Range("B3").Value = "R files created"
$ws2.Range("B18").Value = "R files created"
$ws2.Range("B19").Value = "R files created"

# Selection on this sheet moves from B19 to B3; this sheet stays the active
# (visible) tab.
$ws2.Activate()
$ws2.Range("B3").Select()
